# Remove "Full " from the "Tenant Full Name" header cell, splitting the
# remaining text into two runs ("Tenant " and "Name") the way Word does
# when a formatting touch forces a run break at the edit boundary.

$d = $word.ActiveDocument

# 1) Delete the word "Full " from the header cell text so it reads
#    "Tenant Name".
$delRange = $d.Content
$delRange.Find.Execute("Full ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 2) Re-apply the (unchanged) font color on "Name" to force Word to break
#    it into its own run, matching the original run formatting exactly.
$nameRange = $d.Content
$nameRange.Find.Execute("Name")
$nameRange.Font.Color = 1
$nameRange.Font.Color = 2368548
